$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168982625007629
$ws.Range("B1").Value = 2.287107706069946
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.321753978729248
$ws.Range("E1").Value = 1.229792952537537
